# This workbook contains a weekly price-reporting table for
# "Pepino ensalada" (Vega Central Mapocho de Santiago). A new weekly
# record needs to be inserted above the current row 311, shifting the
# existing rows 311:336 down to 312:337 (the last existing row, 336,
# ends up at 337), and the new row 311 is populated with this week's
# figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 311:336 down by one row so that a new row can be inserted
# at position 311 (xlShiftDown = -4121).
$ws.Rows.Item(311).Insert(-4121)

# Populate the newly inserted row 311 with the new weekly record.
$ws.Range("A311").Value = 9
$ws.Range("B311").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C311").Value = "Metropolitana"
$ws.Range("D311").Value = 44826
$ws.Range("E311").Value = 13
$ws.Range("F311").Value = 100112043
$ws.Range("G311").Value = "Pepino ensalada"
$ws.Range("H311").Value = "Sin especificar"
$ws.Range("I311").Value = "Primera"
$ws.Range("J311").Value = 70
$ws.Range("K311").Value = 22000
$ws.Range("L311").Value = 22000
$ws.Range("M311").Value = 22000
$ws.Range("N311").Value = "$/caja 60 unidades"
$ws.Range("O311").Value = "Región de Arica y Parinacota"
$ws.Range("P311").Value = 367
$ws.Range("Q311").Value = 60
$ws.Range("R311").Value = "Hortaliza"
